# Update cryptos list: price (D) and 1h volume/change (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    # Force the cell to keep its text representation exactly
    # (prevents Excel from auto-converting numeric-looking
    # strings like "1.00" or "11.00" into numbers and losing
    # the trailing zeros / thousand-dot formatting), then
    # restores the original "Normal" style so no formatting
    # change is introduced.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "67.833.16"
$ws.Range("E2").Value = "  +1.76%  "
Set-TextValue $ws.Range("D3") "3.334.49"
$ws.Range("E3").Value = "  +1.87%  "
Set-TextValue $ws.Range("D4") "1.00"
$ws.Range("E4").Value = "  +0.00%  "
Set-TextValue $ws.Range("D5") "582.10"
$ws.Range("E5").Value = "  +1.31%  "
Set-TextValue $ws.Range("D6") "176.29"
$ws.Range("E6").Value = "  +2.20%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +2.27%  "
Set-TextValue $ws.Range("D9") "3.332.30"
$ws.Range("E9").Value = "  +2.06%  "
$ws.Range("E10").Value = "  +7.11%  "
$ws.Range("E11").Value = "  +2.57%  "
$ws.Range("E12").Value = "  +5.42%  "
$ws.Range("E13").Value = "  +3.21%  "
Set-TextValue $ws.Range("D14") "696.02"
$ws.Range("E14").Value = "  +1.00%  "
Set-TextValue $ws.Range("D15") "3.877.18"
$ws.Range("E15").Value = "  +2.13%  "
$ws.Range("E16").Value = "  +2.87%  "
Set-TextValue $ws.Range("D17") "67.806.68"
$ws.Range("E17").Value = "  +1.51%  "
$ws.Range("E18").Value = "  +0.04%  "
Set-TextValue $ws.Range("D19") "3.335.84"
$ws.Range("E19").Value = "  +1.99%  "
$ws.Range("E20").Value = "  +2.06%  "
Set-TextValue $ws.Range("D21") "11.06"
$ws.Range("E21").Value = "  +3.90%  "
Set-TextValue $ws.Range("D22") "0.894"
$ws.Range("E22").Value = "  +1.56%  "
$ws.Range("E23").Value = "  +5.36%  "
Set-TextValue $ws.Range("D24") "17.00"
$ws.Range("E24").Value = "  +1.05%  "
Set-TextValue $ws.Range("D25") "99.93"
$ws.Range("E25").Value = "  +0.48%  "
$ws.Range("E26").Value = "  +2.52%  "
$ws.Range("E27").Value = "  +2.68%  "
Set-TextValue $ws.Range("D28") "9.54"
$ws.Range("E28").Value = "  +5.42%  "
Set-TextValue $ws.Range("D29") "33.03"
$ws.Range("E29").Value = "  -0.64%  "
$ws.Range("E30").Value = "  +3.53%  "
Set-TextValue $ws.Range("D31") "7.04"
$ws.Range("E31").Value = "  +6.53%  "
Set-TextValue $ws.Range("D32") "566.69"
$ws.Range("E32").Value = "  -2.29%  "
Set-TextValue $ws.Range("D33") "11.00"
$ws.Range("E33").Value = "  +2.01%  "
$ws.Range("E34").Value = "  +3.85%  "
Set-TextValue $ws.Range("D35") "57.36"
$ws.Range("E35").Value = "  +4.33%  "
$ws.Range("E36").Value = "  +0.05%  "
Set-TextValue $ws.Range("D37") "3.695.57"
$ws.Range("E37").Value = "  -3.12%  "
Set-TextValue $ws.Range("D38") "3.38"
$ws.Range("E38").Value = "  +2.46%  "
Set-TextValue $ws.Range("D39") "34.83"
$ws.Range("E39").Value = "  +11.44%  "
$ws.Range("E40").Value = "  +5.30%  "
$ws.Range("E41").Value = "  +3.55%  "
$ws.Range("E42").Value = "  +7.47%  "
$ws.Range("E43").Value = "  +2.60%  "
$ws.Range("E44").Value = "  +4.43%  "
Set-TextValue $ws.Range("D45") "3.30"
$ws.Range("E45").Value = "  -1.94%  "
Set-TextValue $ws.Range("D46") "0.0408"
$ws.Range("E46").Value = "  +2.42%  "
Set-TextValue $ws.Range("D47") "2.68"
$ws.Range("E47").Value = "  +6.35%  "
$ws.Range("E48").Value = "  +1.94%  "
$ws.Range("E49").Value = "  -0.04%  "
$ws.Range("E50").Value = "  -0.36%  "
Set-TextValue $ws.Range("D51") "131.41"
$ws.Range("E51").Value = "  +1.83%  "
